$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '26.310.12'
Set-TextValue 'E2' '  +0.82%  '
Set-TextValue 'D3' '1.680.77'
Set-TextValue 'E3' '  +0.87%  '
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '218.36'
Set-TextValue 'E5' '  +0.84%  '
Set-TextValue 'D6' '0.5251'
Set-TextValue 'E6' '  +3.04%  '
Set-TextValue 'E7' '  +0.10%  '
Set-TextValue 'D8' '0.2699'
Set-TextValue 'E8' '  +2.48%  '
Set-TextValue 'D9' '0.06481'
Set-TextValue 'E9' '  +1.22%  '
Set-TextValue 'D10' '21.98'
Set-TextValue 'E10' '  +2.10%  '
Set-TextValue 'D11' '0.07523'
Set-TextValue 'E11' '  +1.56%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.682.56'
Set-TextValue 'E12' '  +0.72%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '4.529'
Set-TextValue 'E13' '  +0.46%  '
Set-TextValue 'D14' '0.5806'
Set-TextValue 'E14' '  +0.21%  '
Set-TextValue 'D15' '0.000008534'
Set-TextValue 'E15' '  -0.05%  '
Set-TextValue 'D16' '64.75'
Set-TextValue 'E16' '  +1.06%  '
Set-TextValue 'D17' '26.329.15'
Set-TextValue 'D18' '4.926'
Set-TextValue 'E18' '  +0.09%  '
Set-TextValue 'E19' '  +0.14%  '
Set-TextValue 'D20' '10.87'
Set-TextValue 'E20' '  +0.68%  '
Set-TextValue 'D21' '189.96'
Set-TextValue 'E21' '  +0.29%  '
Set-TextValue 'D22' '6.206'
Set-TextValue 'E22' '  +0.14%  '
Set-TextValue 'E23' '  +0.04%  '
Set-TextValue 'D24' '145.54'
Set-TextValue 'E24' '  +0.31%  '
Set-TextValue 'D25' '7.811'
Set-TextValue 'E25' '  +2.68%  '
Set-TextValue 'D26' '0.1247'
Set-TextValue 'E26' '  +4.43%  '
Set-TextValue 'D27' '15.78'
Set-TextValue 'E27' '  +1.36%  '
Set-TextValue 'D28' '0.06492'
Set-TextValue 'E28' '  +2.62%  '
Set-TextValue 'D29' '1.358'
Set-TextValue 'E29' '  +4.76%  '
Set-TextValue 'D30' '1.326'
Set-TextValue 'E30' '  +0.78%  '
Set-TextValue 'D31' '3.599'
Set-TextValue 'E31' '  +2.16%  '
Set-TextValue 'D32' '3.596'
Set-TextValue 'E32' '  +2.73%  '
Set-TextValue 'D33' '1.663'
Set-TextValue 'E33' '  +1.82%  '
Set-TextValue 'D34' '1.031'
Set-TextValue 'E34' '  +1.80%  '
Set-TextValue 'D35' '0.6245'
Set-TextValue 'E35' '  +2.85%  '
Set-TextValue 'E36' '  +1.67%  '
Set-TextValue 'E37' '  +2.74%  '
Set-TextValue 'D38' '6.455'
Set-TextValue 'E38' '  +4.77%  '
Set-TextValue 'D39' '1.110.99'
Set-TextValue 'E39' '  +3.37%  '
Set-TextValue 'D40' '0.01623'
Set-TextValue 'E40' '  +1.15%  '
Set-TextValue 'D41' '0.8765'
Set-TextValue 'E41' '  +1.62%  '
Set-TextValue 'E42' '  +0.53%  '
Set-TextValue 'D43' '100.62'
Set-TextValue 'E43' '  -0.36%  '
Set-TextValue 'D44' '1.832.09'
Set-TextValue 'E44' '  +1.00%  '
Set-TextValue 'E45' '  -3.40%  '
Set-TextValue 'D46' '57.03'
Set-TextValue 'E46' '  +1.70%  '
Set-TextValue 'D47' '8.186'
Set-TextValue 'E47' '  +1.67%  '
Set-TextValue 'E48' '  +0.05%  '
Set-TextValue 'D49' '0.05267'
Set-TextValue 'D50' '0.4291'
Set-TextValue 'E50' '  -0.11%  '
Set-TextValue 'E51' '  +2.63%  '
